# Update the team-specific transition-probability matrix on the active sheet.
# The commit "added more games, sped up simulate game logic, and drafted
# optimization logic" re-ran the game simulation with a larger sample, so the
# empirical transition probabilities in rows 2-19 (columns B:S) shift slightly.
# Apply the updated probabilities directly to the affected cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1809815950920245
$ws.Range("C2").Value = 0.588957055214724
$ws.Range("J2").Value = 0.006134969325153374
$ws.Range("O2").Value = 0.003067484662576687
$ws.Range("P2").Value = 0.1226993865030675
$ws.Range("S2").Value = 0.09815950920245399
$ws.Range("B3").Value = 0.01
$ws.Range("C3").Value = 0.03
$ws.Range("J3").Value = 0.02
$ws.Range("P3").Value = 0.76
$ws.Range("S3").Value = 0.18
$ws.Range("J4").Value = 0.02127659574468085
$ws.Range("P4").Value = 0.7659574468085106
$ws.Range("S4").Value = 0.2127659574468085
$ws.Range("B6").Value = 0.04583333333333333
$ws.Range("D6").Value = 0.0125
$ws.Range("J6").Value = 0.2791666666666667
$ws.Range("O6").Value = 0.02083333333333333
$ws.Range("R6").Value = 0.07083333333333333
$ws.Range("S6").Value = 0.3708333333333333
$ws.Range("B7").Value = 0.1302083333333333
$ws.Range("D7").Value = 0.02083333333333333
$ws.Range("F7").Value = 0.046875
$ws.Range("J7").Value = 0.1614583333333333
$ws.Range("O7").Value = 0.02083333333333333
$ws.Range("Q7").Value = 0.171875
$ws.Range("R7").Value = 0.0625
$ws.Range("S7").Value = 0.3854166666666667
$ws.Range("B8").Value = 0.09913043478260869
$ws.Range("D8").Value = 0.02260869565217391
$ws.Range("E8").Value = 0.001739130434782609
$ws.Range("F8").Value = 0.07652173913043478
$ws.Range("J8").Value = 0.1321739130434783
$ws.Range("O8").Value = 0.02086956521739131
$ws.Range("Q8").Value = 0.1878260869565217
$ws.Range("R8").Value = 0.05391304347826087
$ws.Range("S8").Value = 0.4052173913043478
$ws.Range("B9").Value = 0.1233766233766234
$ws.Range("E9").Value = 0.006493506493506494
$ws.Range("F9").Value = 0.07792207792207792
$ws.Range("J9").Value = 0.1298701298701299
$ws.Range("O9").Value = 0.01298701298701299
$ws.Range("Q9").Value = 0.2012987012987013
$ws.Range("R9").Value = 0.06493506493506493
$ws.Range("S9").Value = 0.3831168831168831
$ws.Range("B10").Value = 0.1300411522633745
$ws.Range("D10").Value = 0.02304526748971194
$ws.Range("E10").Value = 0.002469135802469136
$ws.Range("F10").Value = 0.07160493827160494
$ws.Range("J10").Value = 0.1135802469135802
$ws.Range("O10").Value = 0.01975308641975309
$ws.Range("Q10").Value = 0.2197530864197531
$ws.Range("R10").Value = 0.06502057613168724
$ws.Range("S10").Value = 0.3547325102880658
$ws.Range("G11").Value = 0.1360544217687075
$ws.Range("J11").Value = 0.1292517006802721
$ws.Range("K11").Value = 0.1870748299319728
$ws.Range("L11").Value = 0.5238095238095238
$ws.Range("S11").Value = 0.02380952380952381
$ws.Range("G12").Value = 0.7062146892655368
$ws.Range("J12").Value = 0.1638418079096045
$ws.Range("K12").Value = 0.01694915254237288
$ws.Range("L12").Value = 0.03389830508474576
$ws.Range("S12").Value = 0.07909604519774012
$ws.Range("G13").Value = 0.72
$ws.Range("J13").Value = 0.12
$ws.Range("S13").Value = 0.16
$ws.Range("F15").Value = 0.03619909502262444
$ws.Range("H15").Value = 0.2398190045248869
$ws.Range("I15").Value = 0.02714932126696833
$ws.Range("J15").Value = 0.253393665158371
$ws.Range("K15").Value = 0.06334841628959276
$ws.Range("M15").Value = 0.009049773755656109
$ws.Range("O15").Value = 0.06787330316742081
$ws.Range("S15").Value = 0.3031674208144796
$ws.Range("F16").Value = 0.004444444444444444
$ws.Range("H16").Value = 0.1733333333333333
$ws.Range("I16").Value = 0.06222222222222222
$ws.Range("J16").Value = 0.3688888888888889
$ws.Range("K16").Value = 0.1422222222222222
$ws.Range("M16").Value = 0.008888888888888889
$ws.Range("O16").Value = 0.06666666666666667
$ws.Range("S16").Value = 0.1733333333333333
$ws.Range("F17").Value = 0.01458333333333333
$ws.Range("H17").Value = 0.2604166666666667
$ws.Range("I17").Value = 0.07083333333333333
$ws.Range("J17").Value = 0.36875
$ws.Range("K17").Value = 0.1
$ws.Range("M17").Value = 0.025
$ws.Range("O17").Value = 0.05416666666666667
$ws.Range("S17").Value = 0.10625
$ws.Range("F18").Value = 0.04026845637583892
$ws.Range("H18").Value = 0.1476510067114094
$ws.Range("I18").Value = 0.04697986577181208
$ws.Range("J18").Value = 0.3892617449664429
$ws.Range("K18").Value = 0.0738255033557047
$ws.Range("M18").Value = 0.01342281879194631
$ws.Range("O18").Value = 0.1006711409395973
$ws.Range("S18").Value = 0.1879194630872483
$ws.Range("F19").Value = 0.02833078101071975
$ws.Range("H19").Value = 0.2580398162327718
$ws.Range("I19").Value = 0.0666156202143951
$ws.Range("J19").Value = 0.3460949464012251
$ws.Range("K19").Value = 0.09647779479326186
$ws.Range("M19").Value = 0.02297090352220521
$ws.Range("N19").Value = 0.001531393568147014
$ws.Range("O19").Value = 0.06202143950995406
$ws.Range("S19").Value = 0.1179173047473201
